# "Better handle copying paragraph styles"
#
# Each passage block in this template ends with a little run of
# "joiner" paragraphs used by the merge engine:
#   - a (already styled) MSC_Join paragraph
#   - a "[...]" placeholder paragraph      <- was plain/Normal
#   - a blank paragraph right after it     <- was plain/Normal
# The placeholder and the blank paragraph that follows it were not
# getting the MSC_Join paragraph style copied onto them. Fix that by
# giving both of them the MSC_Join style, same as their neighbour.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "[...]`r") {
        $p.Style = "MSC_Join"
        if ($i -lt $count) {
            $next = $d.Paragraphs.Item($i + 1)
            $next.Style = "MSC_Join"
        }
    }
}
